$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the updated values exactly as specified by the diff.
# The "Price" (D) column holds numeric-looking text (e.g. "39.345.33",
# "22.30", "1.00") that must stay as literal text -- a leading apostrophe
# forces Excel to store it as text instead of coercing it to a number
# (which would drop the grouping dots / trailing zeros), matching how the
# source file already stores these values.
$ws.Range("D2").Value = "'39.345.33"
$ws.Range("E2").Value = "  +1.50%  "
$ws.Range("D3").Value = "'2.160.21"
$ws.Range("E3").Value = "  +3.23%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'228.82"
$ws.Range("E5").Value = "  -0.25%  "
$ws.Range("E6").Value = "  +1.08%  "
$ws.Range("D7").Value = "'64.42"
$ws.Range("E7").Value = "  +5.24%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  +3.10%  "
$ws.Range("D10").Value = "'0.0858"
$ws.Range("E10").Value = "  +1.56%  "
$ws.Range("E11").Value = "  -0.22%  "
$ws.Range("D12").Value = "'15.94"
$ws.Range("E12").Value = "  +3.46%  "
$ws.Range("D13").Value = "'2.481.63"
$ws.Range("E13").Value = "  -17.58%  "
$ws.Range("D14").Value = "'22.30"
$ws.Range("E14").Value = "  +0.81%  "
$ws.Range("E15").Value = "  +0.68%  "
$ws.Range("D17").Value = "'2.160.72"
$ws.Range("E17").Value = "  +3.49%  "
$ws.Range("D18").Value = "'39.280.40"
$ws.Range("E18").Value = "  +1.38%  "
$ws.Range("D19").Value = "'71.99"
$ws.Range("E19").Value = "  +0.17%  "
$ws.Range("D20").Value = "'6.15"
$ws.Range("E20").Value = "  +0.93%  "
$ws.Range("E21").Value = "  +1.26%  "
$ws.Range("D22").Value = "'231.23"
$ws.Range("E22").Value = "  +1.55%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").Value = "'2.50"
$ws.Range("E24").Value = "  +5.05%  "
$ws.Range("E25").Value = "  +0.95%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "'172.02"
$ws.Range("E26").Value = "  +0.45%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "'9.52"
$ws.Range("E27").Value = "  -0.38%  "
$ws.Range("E28").Value = "  -1.16%  "
$ws.Range("E29").Value = "  +3.48%  "
$ws.Range("D31").Value = "'2.72"
$ws.Range("E31").Value = "  +9.04%  "
$ws.Range("E32").Value = "  +1.38%  "
$ws.Range("D33").Value = "'4.64"
$ws.Range("E33").Value = "  +2.61%  "
$ws.Range("E34").Value = "  +1.71%  "
$ws.Range("D35").Value = "'7.09"
$ws.Range("E35").Value = "  +7.97%  "
$ws.Range("E36").Value = "  +1.43%  "
$ws.Range("E37").Value = "  +1.26%  "
$ws.Range("D38").Value = "'3.59"
$ws.Range("E38").Value = "  -0.37%  "
$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "  -0.11%  "
$ws.Range("D40").Value = "'104.16"
$ws.Range("E40").Value = "  +2.99%  "
$ws.Range("E41").Value = "  +0.09%  "
$ws.Range("D42").Value = "'17.79"
$ws.Range("E42").Value = "  -1.29%  "
$ws.Range("D43").Value = "'1.540.91"
$ws.Range("E43").Value = "  +0.52%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Value = "'1.18"
$ws.Range("E44").Value = "  +3.58%  "
$ws.Range("B45").Value = "FTXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D45").Value = "'4.32"
$ws.Range("E45").Value = "  +4.75%  "
$ws.Range("E46").Value = "  +0.62%  "
$ws.Range("D47").Value = "'0.0925"
$ws.Range("E47").Value = "  +0.92%  "
$ws.Range("E48").Value = "  +6.01%  "
$ws.Range("D49").Value = "'7.80"
$ws.Range("E49").Value = "  +1.84%  "
$ws.Range("D50").Value = "'2.364.69"
$ws.Range("E50").Value = "  +3.37%  "
$ws.Range("D51").Value = "'2.96"
$ws.Range("E51").Value = "  -0.33%  "
